$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.710.98"
$ws.Range("E2").Value = "  +2.07%  "
$ws.Range("D3").Value = "'2.042.29"
$ws.Range("E3").Value = "  +1.56%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'230.00"
$ws.Range("E5").Value = "  +1.99%  "
$ws.Range("D6").Value = "'0.608"
$ws.Range("E6").Value = "  +0.20%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'56.33"
$ws.Range("E8").Value = "  +3.42%  "
$ws.Range("D9").Value = "'0.383"
$ws.Range("E9").Value = "  +1.29%  "
$ws.Range("D10").Value = "'0.0803"
$ws.Range("E10").Value = "  +2.61%  "
$ws.Range("E11").Value = "  -0.70%  "
$ws.Range("D12").Value = "'2.342.60"
$ws.Range("E12").Value = "  +1.34%  "
$ws.Range("D13").Value = "'14.43"
$ws.Range("E13").Value = "  +1.31%  "
$ws.Range("D14").Value = "'20.31"
$ws.Range("E14").Value = "  +0.18%  "
$ws.Range("D15").Value = "'5.23"
$ws.Range("E15").Value = "  +2.13%  "
$ws.Range("D17").Value = "'2.038.59"
$ws.Range("E17").Value = "  +1.44%  "
$ws.Range("D18").Value = "'37.616.40"
$ws.Range("E18").Value = "  +2.05%  "
$ws.Range("E19").Value = "  +0.27%  "
$ws.Range("D20").Value = "'69.21"
$ws.Range("E20").Value = "  +0.77%  "
$ws.Range("E21").Value = "  +1.16%  "
$ws.Range("D22").Value = "'223.71"
$ws.Range("E22").Value = "  -0.82%  "
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("E24").Value = "  +1.52%  "
$ws.Range("E25").Value = "  +3.38%  "
$ws.Range("D26").Value = "'165.40"
$ws.Range("E26").Value = "  +0.32%  "
$ws.Range("D27").Value = "'9.19"
$ws.Range("E27").Value = "  -0.19%  "
$ws.Range("E28").Value = "  +5.85%  "
$ws.Range("E29").Value = "  +2.02%  "
$ws.Range("E30").Value = "  +0.16%  "
$ws.Range("E31").Value = "  +1.06%  "
$ws.Range("E32").Value = "  -0.17%  "
$ws.Range("E33").Value = "  -0.57%  "
$ws.Range("E34").Value = "  +2.36%  "
$ws.Range("D35").Value = "'2.02"
$ws.Range("E35").Value = "  +9.04%  "
$ws.Range("D36").Value = "'2.35"
$ws.Range("E36").Value = "  +0.16%  "
$ws.Range("D37").Value = "'5.92"
$ws.Range("E37").Value = "  +10.64%  "
$ws.Range("D38").Value = "'3.26"
$ws.Range("E38").Value = "  +4.06%  "
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").Value = "'97.63"
$ws.Range("E40").Value = "  +2.72%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "'1.478.29"
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("B42").Value = "Cronos"
$ws.Range("C42").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D42").Value = "'0.0948"
$ws.Range("E42").Value = "  +3.56%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "'0.0216"
$ws.Range("E43").Value = "  -0.57%  "
$ws.Range("B44").Value = "HuobiToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D44").Value = "'2.85"
$ws.Range("E44").Value = "  +2.56%  "
$ws.Range("D45").Value = "'16.62"
$ws.Range("E45").Value = "  -0.85%  "
$ws.Range("D46").Value = "'4.16"
$ws.Range("E46").Value = "  +16.84%  "
$ws.Range("E47").Value = "  -1.27%  "
$ws.Range("E48").Value = "  +1.19%  "
$ws.Range("E49").Value = "  -0.87%  "
$ws.Range("E50").Value = "  +0.50%  "
$ws.Range("D51").Value = "'2.230.93"
$ws.Range("E51").Value = "  +1.49%  "
